$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3123
$ws1.Range("F5").Value = 26
$ws1.Range("F8").Value = 315
$ws1.Range("F9").Value = 7325
$ws1.Range("F10").Value = 68
$ws1.Range("F11").Value = 147
$ws1.Range("F12").Value = 57
$ws1.Range("F13").Value = 30
$ws1.Range("F14").Value = 427
$ws1.Range("F15").Value = 84
$ws1.Range("F16").Value = 84
$ws1.Range("F17").Value = 1946
$ws1.Range("F18").Value = 1788
$ws1.Range("F19").Value = 1079
$ws1.Range("F21").Value = 70
$ws1.Range("F22").Value = 1820
$ws1.Range("F23").Value = 1370
$ws1.Range("F24").Value = 1239
$ws1.Range("F25").Value = 642
$ws1.Range("F26").Value = 51
$ws1.Range("F27").Value = 1123
$ws1.Range("F28").Value = 119
$ws1.Range("F29").Value = 533
$ws1.Range("F32").Value = 2692
$ws1.Range("F33").Value = 1515
$ws1.Range("F34").Value = 3201
$ws1.Range("F35").Value = 2183
$ws1.Range("F36").Value = 148
$ws1.Range("F37").Value = 220
$ws1.Range("F40").Value = 44
$ws1.Range("F42").Value = 375
$ws1.Range("F43").Value = 154
$ws1.Range("F44").Value = 517
$ws1.Range("F45").Value = 250
$ws1.Range("F47").Value = 765
$ws1.Range("F48").Value = 434
$ws1.Range("F49").Value = 1
$ws1.Range("F50").Value = 122

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("C24").Value = "上海·触手猴动漫钢琴音乐演奏会  Marasy Piano Live Asia Tour Prelive "
$ws2.Range("F24").Value = 82

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 1842
$ws3.Range("F8").Value = 2893
$ws3.Range("F10").Value = 1116
$ws3.Range("F12").Value = 420
$ws3.Range("F13").Value = 1825
$ws3.Range("F14").Value = 8156

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3123
$ws4.Range("F6").Value = 1842
$ws4.Range("F7").Value = 2893
$ws4.Range("F9").Value = 1116
$ws4.Range("F10").Value = 68
$ws4.Range("F11").Value = 420
$ws4.Range("F12").Value = 57
$ws4.Range("F13").Value = 30
$ws4.Range("F14").Value = 427
$ws4.Range("F15").Value = 84
$ws4.Range("F16").Value = 84
$ws4.Range("F18").Value = 1079
$ws4.Range("F20").Value = 70
$ws4.Range("F21").Value = 1820
$ws4.Range("F22").Value = 1370
$ws4.Range("F23").Value = 1239
$ws4.Range("F24").Value = 642
$ws4.Range("F25").Value = 51
$ws4.Range("F26").Value = 1123
$ws4.Range("F28").Value = 119
$ws4.Range("F31").Value = 533
$ws4.Range("F34").Value = 2692
$ws4.Range("F35").Value = 1515
$ws4.Range("F36").Value = 3209
$ws4.Range("F37").Value = 2183
$ws4.Range("F38").Value = 148
$ws4.Range("F39").Value = 220
$ws4.Range("F43").Value = 154
$ws4.Range("C44").Value = "上海·触手猴动漫钢琴音乐演奏会  Marasy Piano Live Asia Tour Prelive "
$ws4.Range("F44").Value = 82
$ws4.Range("F45").Value = 517
$ws4.Range("F46").Value = 250
$ws4.Range("F48").Value = 434
